# Edit script for undyingkingdoms modifiers.xlsx
# - Renames "Dwarven Made" -> "Dwarven Steel" (row 2)
# - Adjusts Goblin Expendable birth rate (10% -> 15%) and Infighting happiness (-2 -> -1)
# - Inserts a new Goblin "Sneaky" infiltration modifier (row 8), shifting Human and all
#   subsequent background modifiers down by two rows, and renames/retunes several of them
# - Extends the sheet with blank rows through row 45 to fix mobile compare-table formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Dwarf / Dwarven Made -> Dwarven Steel ---
$ws.Range("D2").Value = "Dwarven Steel"
$ws.Range("G2").Value = "Dwarven Steel: +15% Production Per Worker"

# --- Row 6: Goblin / Expendable birth rate 10% -> 15% ---
$ws.Range("F6").Value = 0.15
$ws.Range("G6").Value = "Expendable: +15% Birth Rate"

# --- Row 7: Goblin / Infighting happiness -2 -> -1 ---
$ws.Range("F7").Value = -1
$ws.Range("G7").Value = "Infighting: -1 Happiness"

# --- Rows 8-15: re-sequence background/race modifier rows to make room for the new
#     "Sneaky" Goblin row, and retune/rename several Background modifiers ---
$ws.Range("B8").Value = "Goblin"
$ws.Range("C8").Value = "'"
$ws.Range("D8").Value = "Sneaky"
$ws.Range("E8").Value = "Infiltration Success Modifier"
$ws.Range("F8").Value = 0.1
$ws.Range("G8").Value = "Sneaky: +10% Infiltration Success"

$ws.Range("B9").Value = "Human"
$ws.Range("C9").Value = "'"
$ws.Range("D9").Value = "'"
$ws.Range("E9").Value = "'"
$ws.Range("F9").Value = "'"
$ws.Range("G9").Value = "'"

$ws.Range("B10").Value = "'"
$ws.Range("C10").Value = "Engineer"
$ws.Range("D10").Value = "Artisan"
$ws.Range("E10").Value = "Production Per Worker Modifier"
$ws.Range("F10").Value = 0.2
$ws.Range("G10").Value = "Artisan: +20% Production Per Worker"

$ws.Range("B11").Value = "'"
$ws.Range("C11").Value = "Engineer"
$ws.Range("D11").Value = "Craftsman"
$ws.Range("E11").Value = "Buildings Built Per Day Modifier"
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = "Craftsman: +2 Buildings Built Per Day"

$ws.Range("B12").Value = "'"
$ws.Range("C12").Value = "Merchant"
$ws.Range("D12").Value = "Silver Tongue"
$ws.Range("E12").Value = "Income Modifier"
$ws.Range("F12").Value = 0.15
$ws.Range("G12").Value = "Silver Tongue: +15% Income"

$ws.Range("B13").Value = "'"
$ws.Range("C13").Value = "Rogue"
$ws.Range("D13").Value = "Army of Shadows"
$ws.Range("E13").Value = "Amount Of Thieves Modifier"
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = "Army of Shadows: +2 Amount Of Thieves"

$ws.Range("B14").Value = "'"
$ws.Range("C14").Value = "Rogue"
$ws.Range("D14").Value = "Master of Disguise"
$ws.Range("E14").Value = "Infiltration Success Modifier"
$ws.Range("F14").Value = 0.15
$ws.Range("G14").Value = "Master of Disguise: +15% Infiltration Success"

$ws.Range("B15").Value = "'"
$ws.Range("C15").Value = "Warlord"
$ws.Range("D15").Value = "Relentless"
$ws.Range("E15").Value = "Offensive Power Modifier"
$ws.Range("F15").Value = 0.15
$ws.Range("G15").Value = "Relentless: +15% Offensive Power"

# Rows 1-13 already carry the bold/centered index-column style; rows 14+ are brand
# new, so their "#" index column (A) needs both the value and the matching style.
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13

# --- Rows 16-45: extend sheet with blank index rows (fixes mobile compare table) ---
for ($r = 16; $r -le 45; $r++) {
    $ws.Range("A" + $r).Value = $r - 2
    $ws.Range("B" + $r + ":G" + $r).Value = "'"
}

# Strip the quote-prefix / stray style introduced by the "'" placeholders above so the
# new blank cells end up as plain empty-string cells (like the rest of the sheet) --
# "Z1" is outside the used range and was never written to, so it carries the default,
# un-styled format to paste over every blank cell touched by this script.
$ws.Range("Z1").Copy()
$ws.Range("B8:G15").PasteSpecial(-4122)
$ws.Range("B16:G45").PasteSpecial(-4122)

# Re-apply the bold/centered/thin-border "#" index-column style (lost by the blanket
# paste above) to the new index rows.
$ws.Range("A2").Copy()
$ws.Range("A14:A45").PasteSpecial(-4122)

$excel.CutCopyMode = 0

